$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row: "_old" -> "_FV2404" and "_new" -> "_FV2410" ---
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")
$labels  = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $labels.Count; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $labels[$i] + "_FV2404"
    $ws.Range($newCols[$i] + "1").Value = $labels[$i] + "_FV2410"
}

# --- 2. Turn the used range A1:U57 into an Excel Table ("Table1") ---
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U57"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header renamed, Table1 created, header row frozen."
